# Update Cxcl5-Cxcr2 LR-pair sheet with new TPM-derived values.
#
# The old data had 4 data rows (ECs->ECs, ECs->FAPs, FAPs->ECs, FAPs->FAPs).
# The refreshed export only keeps the FAPs-sourced rows (FAPs->ECs, FAPs->FAPs)
# recomputed against the new TPM numbers, so we drop the two stale
# ECs-sourced rows and rewrite the remaining two with the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete "ECs" sending-cluster rows (old rows 2 and 3).
$ws.Rows("2:3").Delete()

# Row 2: FAPs -> Cxcl5 -> Cxcr2 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl5"
$ws.Range("C2").Value = "Cxcr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05862133333333333
$ws.Range("H2").Value = 0.175864
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01650666666666667
$ws.Range("N2").Value = 0.04952
$ws.Range("O2").Value = 0.795859985214233
$ws.Range("P2").Value = 0.795859985214233
$ws.Range("Q2").Value = 0.0009676428088888888
$ws.Range("R2").Value = 0.008708785280000001
$ws.Range("S2").Value = 0.795859985214233
$ws.Range("T2").Value = 0.795859985214233

# Row 3: FAPs -> Cxcl5 -> Cxcr2 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl5"
$ws.Range("C3").Value = "Cxcr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05862133333333333
$ws.Range("H3").Value = 0.175864
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.004234
$ws.Range("N3").Value = 0.012702
$ws.Range("O3").Value = 0.2041400147857671
$ws.Range("P3").Value = 0.2041400147857671
$ws.Range("Q3").Value = 0.0002482027253333333
$ws.Range("R3").Value = 0.002233824528
$ws.Range("S3").Value = 0.2041400147857671
$ws.Range("T3").Value = 0.2041400147857671
